$d = $word.ActiveDocument

# Replace oldText -> newText, scoping the Find to the exact span where oldText
# occurs so we never touch unrelated text elsewhere in the document.
function Replace-Scoped($oldText, $newText) {
    $text = $d.Content.Text
    $idx = $text.IndexOf($oldText)
    if ($idx -ge 0) {
        $r = $d.Range($idx, $idx + $oldText.Length)
        $r.Find.Execute($oldText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newText, 2)
    }
}

# 1. Créditos-trabalho: 1 -> 0
Replace-Scoped "Créditos-trabalho: 1" "Créditos-trabalho: 0"

# 2. Carga horária: 90 h -> 60 h
Replace-Scoped "Carga horária: 90 h" "Carga horária: 60 h"

# 3. Ativação: 01/01/2020 -> 01/01/2025
Replace-Scoped "Ativação: 01/01/2020" "Ativação: 01/01/2025"

# 4. Append sentence to Portuguese "Programa" paragraph
Replace-Scoped "Infraestrutura de dados espaciais. Cartografia digital." "Infraestrutura de dados espaciais. Cartografia digital. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# 5. Append sentence to English "Programa" paragraph
Replace-Scoped "Spatial data infrastructure. Digital cartography." "Spatial data infrastructure. Digital cartography. The discipline may have didactic trips to complement the content of the discipline."
